$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------------
# Shape "Content Placeholder 2" (id=3) - left-hand agenda box
# ---------------------------------------------------------------------------
$sh1 = $s.Shapes.Item(2)
$tr1 = $sh1.TextFrame.TextRange

# Paragraph 1: "IG: Monday March 15 (3h)" + " " (2 runs) -> single run with
# trailing space, still bold. Force a text reset so the two runs collapse
# into one before writing the final text.
$pa1 = $tr1.Paragraphs(1, 1)
$pa1.Text = "X"
$pa1b = $tr1.Paragraphs(1, 1)
$pa1b.Text = "IG: Monday March 15 (3h) "
$pa1b.Font.Bold = -1

# Paragraph 2: "Joint Sessions " (lvl 1) becomes bold.
$pa2 = $tr1.Paragraphs(2, 1)
$pa2.Font.Bold = -1

# Paragraph 3: "WG: Wednesday March 17 (3h)" loses its bold.
$pa3 = $tr1.Paragraphs(3, 1)
$pa3.Font.Bold = 0

# Paragraph 5: "IG: Thursday March 18 (2h)" loses its bold.
$pa5 = $tr1.Paragraphs(5, 1)
$pa5.Font.Bold = 0

# ---------------------------------------------------------------------------
# Shape "Content Placeholder 2" (id=8) - right-hand agenda textbox
# ---------------------------------------------------------------------------
$sh2 = $s.Shapes.Item(6)
$tr2 = $sh2.TextFrame.TextRange

# Paragraph 1: "WG: Monday March 22" + " " + "(3h)" (3 runs) -> single run,
# no longer bold.
$qa1 = $tr2.Paragraphs(1, 1)
$qa1.Text = "X"
$qa1b = $tr2.Paragraphs(1, 1)
$qa1b.Text = "WG: Monday March 22 (3h)"
$qa1b.Font.Bold = 0

# Paragraph 4: "WG: Wednesday March 24 (3h)" loses its bold.
$qa4 = $tr2.Paragraphs(4, 1)
$qa4.Font.Bold = 0

# Paragraph 7: "WG: Thursday March 25 (2h)" loses its bold.
$qa7 = $tr2.Paragraphs(7, 1)
$qa7.Font.Bold = 0
